## Dependencies: Remove lemmalist-greek; Utils: Remove lemmalist-greek's Greek (Ancient) lemma list
##
## Deletes the worksheet row describing the "lemmalist-greek" dependency
## (entire row, shifting everything below it up by one), then rebuilds the
## hyperlinks collection (which this runtime does not auto-shift when rows
## are deleted) from the literal URL text already present in columns B
## (Home Page) and F (License URL) of every remaining data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Locate and delete the "lemmalist-greek" row -------------------------
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$targetRow = -1
for ($r = 1; $r -le $lastRow; $r++) {
    $nameVal = $ws.Cells.Item($r, 1).Value2
    if ($nameVal -eq "lemmalist-greek") {
        $targetRow = $r
        break
    }
}

if ($targetRow -gt 0) {
    $ws.Rows($targetRow).Delete()
}

# --- 2. Rebuild hyperlinks ---------------------------------------------------
# Row deletion does not shift/clean the <hyperlinks> table in this runtime,
# so drop every hyperlink and recreate them from the (already-shifted) cell
# text, which stores "Address" or "Address#SubAddress".
$ws.Hyperlinks.Delete()

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

function Add-RowHyperlink($sheet, $rowNum, $colLetter, $colIndex) {
    $cell = $sheet.Cells.Item($rowNum, $colIndex)
    $text = $cell.Value2
    if ([string]::IsNullOrEmpty($text)) {
        return
    }
    $hashPos = $text.IndexOf("#")
    if ($hashPos -ge 0) {
        $address = $text.Substring(0, $hashPos)
        $subAddress = $text.Substring($hashPos + 1)
    } else {
        $address = $text
        $subAddress = ""
    }
    $ref = "$colLetter$rowNum"
    $sheet.Hyperlinks.Add($sheet.Range($ref), $address, $subAddress)
}

for ($r = 2; $r -le $lastRow; $r++) {
    Add-RowHyperlink $ws $r "B" 2
    Add-RowHyperlink $ws $r "F" 6
}

# --- 3. Restore the saved selection state -----------------------------------
$ws.Range("D11").Select()
